# Sort the data table (A1:D61, header in row 1) ascending by the ID
# column (A), then fix up the last row: after sorting, the row that
# carries ID 60 still shows the old Name/Country/Product values that
# belonged there before the sort - overwrite them with the values from
# the row directly above (ID 59) to match the corrected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A61"))
$ws.Sort.SetRange($ws.Range("A1:D61"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

$ws.Range("B61").Value = $ws.Range("B60").Value2
$ws.Range("C61").Value = $ws.Range("C60").Value2
$ws.Range("D61").Value = $ws.Range("D60").Value2
